# This script refreshes Market Board price/profit columns (H:N) across
# several Leve-tracking worksheets, matching an upstream scheduled-runner
# price sync. Each block targets one worksheet; within it, cells are
# updated row by row for the leve rows whose cached prices changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Range("H132").Value = 18420.34
$ws.Range("I132").Value = 19596.254
$ws.Range("J132").Value = 2251.5
$ws.Range("K132").Value = 58788.762
$ws.Range("L132").Value = 6754.5
$ws.Range("M132").Value = -56258.762
$ws.Range("N132").Value = -11814.5
# Row 134
$ws.Range("H134").Value = 40870
$ws.Range("J134").Value = 41022.855
$ws.Range("L134").Value = 41022.855
$ws.Range("N134").Value = -51162.855
# Row 135
$ws.Range("H135").Value = 644.2632
$ws.Range("I135").Value = 626.5405
$ws.Range("J135").Value = 1300
$ws.Range("K135").Value = 5638.8645
$ws.Range("L135").Value = 11700
$ws.Range("M135").Value = -3103.8645
$ws.Range("N135").Value = -16770
# Row 138
$ws.Range("H138").Value = 2876.74
$ws.Range("I138").Value = 1373.3684
$ws.Range("J138").Value = 4869.5815
$ws.Range("K138").Value = 4120.1052
$ws.Range("L138").Value = 14608.7445
$ws.Range("M138").Value = 1019.8948
$ws.Range("N138").Value = -24888.7445
# Row 141
$ws.Range("H141").Value = 4233.1577
$ws.Range("I141").Value = 1541.2122
$ws.Range("J141").Value = 22000
$ws.Range("K141").Value = 4623.6366
$ws.Range("L141").Value = 66000
$ws.Range("M141").Value = 556.3634000000002
$ws.Range("N141").Value = -76360

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2634.02
$ws.Range("I32").Value = 2635.885
$ws.Range("J32").Value = 2621.5386
$ws.Range("K32").Value = 2635.885
$ws.Range("L32").Value = 2621.5386
$ws.Range("M32").Value = -2348.885
$ws.Range("N32").Value = -3195.5386
# Row 94
$ws.Range("H94").Value = 18000
$ws.Range("J94").Value = 18000
$ws.Range("L94").Value = 18000
$ws.Range("N94").Value = -19802
# Row 102
$ws.Range("H102").Value = 1566.6666
$ws.Range("I102").Value = 1566.6666
$ws.Range("K102").Value = 1566.6666
$ws.Range("M102").Value = 55.33339999999998
# Row 109
$ws.Range("H109").Value = 15915.4
$ws.Range("J109").Value = 15915.4
$ws.Range("L109").Value = 15915.4
$ws.Range("N109").Value = -18689.4
# Row 110
$ws.Range("H110").Value = 15167.543
$ws.Range("I110").Value = 16683.355
$ws.Range("J110").Value = 3420
$ws.Range("K110").Value = 16683.355
$ws.Range("L110").Value = 3420
$ws.Range("M110").Value = -14638.355
$ws.Range("N110").Value = -7510

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2130.8572
$ws.Range("I134").Value = 1105.8948
$ws.Range("J134").Value = 3348
$ws.Range("K134").Value = 3317.6844
$ws.Range("L134").Value = 10044
$ws.Range("M134").Value = -782.6844000000001
$ws.Range("N134").Value = -15114

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1178.7925
$ws.Range("I58").Value = 886.70215
$ws.Range("J58").Value = 3466.8333
$ws.Range("K58").Value = 886.70215
$ws.Range("L58").Value = 3466.8333
$ws.Range("M58").Value = -683.70215
$ws.Range("N58").Value = -3872.8333
# Row 132
$ws.Range("H132").Value = 1992.7073
$ws.Range("I132").Value = 1439.4062
$ws.Range("J132").Value = 3960
$ws.Range("K132").Value = 4318.2186
$ws.Range("L132").Value = 11880
$ws.Range("M132").Value = -1788.2186
$ws.Range("N132").Value = -16940
# Row 134
$ws.Range("H134").Value = 1705.6735
$ws.Range("I134").Value = 1796.909
$ws.Range("J134").Value = 1517.5
$ws.Range("K134").Value = 5390.727000000001
$ws.Range("L134").Value = 4552.5
$ws.Range("M134").Value = -2855.727000000001
$ws.Range("N134").Value = -9622.5
# Row 136
$ws.Range("H136").Value = 1178.7925
$ws.Range("I136").Value = 886.70215
$ws.Range("J136").Value = 3466.8333
$ws.Range("K136").Value = 2660.10645
$ws.Range("L136").Value = 10400.4999
$ws.Range("M136").Value = -110.1064499999998
$ws.Range("N136").Value = -15500.4999

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 549164.4
$ws.Range("I5").Value = 517.3913
$ws.Range("J5").Value = 1951262.4
$ws.Range("K5").Value = 1552.1739
$ws.Range("L5").Value = 5853787.199999999
$ws.Range("M5").Value = -1440.1739
$ws.Range("N5").Value = -5854011.199999999
# Row 122
$ws.Range("H122").Value = 1274.5862
$ws.Range("I122").Value = 442.125
$ws.Range("J122").Value = 2299.1538
$ws.Range("K122").Value = 3979.125
$ws.Range("L122").Value = 20692.3842
$ws.Range("M122").Value = -1529.125
$ws.Range("N122").Value = -25592.3842
# Row 129
$ws.Range("H129").Value = 1787.3
$ws.Range("I129").Value = 743.3333
$ws.Range("J129").Value = 2483.2778
$ws.Range("K129").Value = 2229.9999
$ws.Range("L129").Value = 7449.8334
$ws.Range("M129").Value = 2770.0001
$ws.Range("N129").Value = -17449.8334
# Row 131
$ws.Range("H131").Value = 3743.72
$ws.Range("I131").Value = 422.7647
$ws.Range("J131").Value = 5454.515
$ws.Range("K131").Value = 1268.2941
$ws.Range("L131").Value = 16363.545
$ws.Range("M131").Value = 3771.7059
$ws.Range("N131").Value = -26443.545
# Row 135
$ws.Range("H135").Value = 549164.4
$ws.Range("I135").Value = 517.3913
$ws.Range("J135").Value = 1951262.4
$ws.Range("K135").Value = 4656.5217
$ws.Range("L135").Value = 17561361.6
$ws.Range("M135").Value = -2121.5217
$ws.Range("N135").Value = -17566431.6

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4926.426
$ws.Range("I70").Value = 4365.7188
$ws.Range("J70").Value = 5742
$ws.Range("K70").Value = 4365.7188
$ws.Range("L70").Value = 5742
$ws.Range("M70").Value = -4095.7188
$ws.Range("N70").Value = -6282
# Row 73
$ws.Range("H73").Value = 4926.426
$ws.Range("I73").Value = 4365.7188
$ws.Range("J73").Value = 5742
$ws.Range("K73").Value = 4365.7188
$ws.Range("L73").Value = 5742
$ws.Range("M73").Value = -3429.7188
$ws.Range("N73").Value = -7614

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 783.3
$ws.Range("I61").Value = 805.375
$ws.Range("J61").Value = 695
$ws.Range("K61").Value = 805.375
$ws.Range("L61").Value = 695
$ws.Range("M61").Value = -603.375
$ws.Range("N61").Value = -1099
# Row 113
$ws.Range("H113").Value = 783.3
$ws.Range("I113").Value = 805.375
$ws.Range("J113").Value = 695
$ws.Range("K113").Value = 805.375
$ws.Range("L113").Value = 695
$ws.Range("M113").Value = 1364.625
$ws.Range("N113").Value = -5035
# Row 132
$ws.Range("H132").Value = 7257.2
$ws.Range("I132").Value = 8237.166999999999
$ws.Range("J132").Value = 5400.421
$ws.Range("K132").Value = 24711.501
$ws.Range("L132").Value = 16201.263
$ws.Range("M132").Value = -22181.501
$ws.Range("N132").Value = -21261.263
# Row 135
$ws.Range("H135").Value = 50167.867
$ws.Range("J135").Value = 50167.867
$ws.Range("L135").Value = 50167.867
$ws.Range("N135").Value = -60307.867

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 40067.23
$ws.Range("I122").Value = 63975.062
$ws.Range("J122").Value = 1814.7
$ws.Range("K122").Value = 191925.186
$ws.Range("L122").Value = 5444.1
$ws.Range("M122").Value = -189475.186
$ws.Range("N122").Value = -10344.1
# Row 132
$ws.Range("H132").Value = 1466.8254
$ws.Range("I132").Value = 1306.5918
$ws.Range("J132").Value = 2027.6428
$ws.Range("K132").Value = 3919.7754
$ws.Range("L132").Value = 6082.928400000001
$ws.Range("M132").Value = -1389.7754
$ws.Range("N132").Value = -11142.9284
